$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2768.3333
$ws.Range("I34").Value = 2768.3333
$ws.Range("K34").Value = 2768.3333
$ws.Range("M34").Value = -2565.3333
$ws.Range("H36").Value = 2768.3333
$ws.Range("I36").Value = 2768.3333
$ws.Range("K36").Value = 2768.3333
$ws.Range("M36").Value = -2053.3333
$ws.Range("H80").Value = 831.51514
$ws.Range("J80").Value = 801.7826
$ws.Range("L80").Value = 2405.3478
$ws.Range("N80").Value = -4401.3478
$ws.Range("H83").Value = 831.51514
$ws.Range("J83").Value = 801.7826
$ws.Range("L83").Value = 7216.0434
$ws.Range("N83").Value = -17200.0434
$ws.Range("H103").Value = 993.88
$ws.Range("I103").Value = 786.8125
$ws.Range("J103").Value = 1362
$ws.Range("K103").Value = 2360.4375
$ws.Range("L103").Value = 4086
$ws.Range("M103").Value = -1774.4375
$ws.Range("N103").Value = -5258
$ws.Range("H109").Value = 65406.617
$ws.Range("J109").Value = 65406.617
$ws.Range("L109").Value = 65406.617
$ws.Range("N109").Value = -68180.617
$ws.Range("H111").Value = 2935.4443
$ws.Range("I111").Value = 2921.4375
$ws.Range("K111").Value = 8764.3125
$ws.Range("M111").Value = -5697.3125
$ws.Range("H132").Value = 1624.921
$ws.Range("I132").Value = 1074.1111
$ws.Range("J132").Value = 2976.9092
$ws.Range("K132").Value = 3222.3333
$ws.Range("L132").Value = 8930.7276
$ws.Range("M132").Value = -692.3333000000002
$ws.Range("N132").Value = -13990.7276
$ws.Range("H138").Value = 2449.5186
$ws.Range("J138").Value = 3325.5
$ws.Range("L138").Value = 9976.5
$ws.Range("N138").Value = -20256.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5037.4
$ws.Range("I61").Value = 5327
$ws.Range("J61").Value = 3396.3333
$ws.Range("K61").Value = 5327
$ws.Range("L61").Value = 3396.3333
$ws.Range("M61").Value = -5115
$ws.Range("N61").Value = -3820.3333
$ws.Range("H74").Value = 1608.4333
$ws.Range("J74").Value = 1767.2222
$ws.Range("L74").Value = 1767.2222
$ws.Range("N74").Value = -3515.2222
$ws.Range("H77").Value = 1608.4333
$ws.Range("J77").Value = 1767.2222
$ws.Range("L77").Value = 8836.110999999999
$ws.Range("N77").Value = -17572.111
$ws.Range("H132").Value = 2063.4285
$ws.Range("I132").Value = 2101.6843
$ws.Range("K132").Value = 6305.0529
$ws.Range("M132").Value = -3775.0529
$ws.Range("H136").Value = 5037.4
$ws.Range("I136").Value = 5327
$ws.Range("J136").Value = 3396.3333
$ws.Range("K136").Value = 15981
$ws.Range("L136").Value = 10188.9999
$ws.Range("M136").Value = -13431
$ws.Range("N136").Value = -15288.9999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1686.1818
$ws.Range("I20").Value = 1170
$ws.Range("J20").Value = 2116.3333
$ws.Range("K20").Value = 1170
$ws.Range("L20").Value = 2116.3333
$ws.Range("M20").Value = -923
$ws.Range("N20").Value = -2610.3333
$ws.Range("H99").Value = 2626
$ws.Range("I99").Value = 1173.5
$ws.Range("K99").Value = 1173.5
$ws.Range("M99").Value = 324.5
$ws.Range("H134").Value = 4334.3
$ws.Range("I134").Value = 5257.1665
$ws.Range("J134").Value = 2950
$ws.Range("K134").Value = 15771.4995
$ws.Range("L134").Value = 8850
$ws.Range("M134").Value = -13236.4995
$ws.Range("N134").Value = -13920

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4038.4783
$ws.Range("I31").Value = 1622.8462
$ws.Range("J31").Value = 7178.8
$ws.Range("K31").Value = 1622.8462
$ws.Range("L31").Value = 7178.8
$ws.Range("M31").Value = -1327.8462
$ws.Range("N31").Value = -7768.8
$ws.Range("H34").Value = 4038.4783
$ws.Range("I34").Value = 1622.8462
$ws.Range("J34").Value = 7178.8
$ws.Range("K34").Value = 1622.8462
$ws.Range("L34").Value = 7178.8
$ws.Range("M34").Value = -1420.8462
$ws.Range("N34").Value = -7582.8
$ws.Range("H39").Value = 5141.5
$ws.Range("I39").Value = 5141.5
$ws.Range("K39").Value = 5141.5
$ws.Range("M39").Value = -4750.5
$ws.Range("H49").Value = 5141.5
$ws.Range("I49").Value = 5141.5
$ws.Range("K49").Value = 5141.5
$ws.Range("M49").Value = -4959.5
$ws.Range("H63").Value = 38333.332
$ws.Range("J63").Value = 38333.332
$ws.Range("L63").Value = 38333.332
$ws.Range("N63").Value = -39705.332
$ws.Range("H64").Value = 42623
$ws.Range("I64").Value = 10246
$ws.Range("K64").Value = 10246
$ws.Range("M64").Value = -9998
$ws.Range("H66").Value = 38333.332
$ws.Range("J66").Value = 38333.332
$ws.Range("L66").Value = 114999.996
$ws.Range("N66").Value = -121863.996
$ws.Range("H67").Value = 42623
$ws.Range("I67").Value = 10246
$ws.Range("K67").Value = 10246
$ws.Range("M67").Value = -9388
$ws.Range("H69").Value = 31245
$ws.Range("I69").Value = 33326.668
$ws.Range("J69").Value = 25000
$ws.Range("K69").Value = 33326.668
$ws.Range("L69").Value = 25000
$ws.Range("M69").Value = -32577.668
$ws.Range("N69").Value = -26498
$ws.Range("H72").Value = 31245
$ws.Range("I72").Value = 33326.668
$ws.Range("J72").Value = 25000
$ws.Range("K72").Value = 99980.00399999999
$ws.Range("L72").Value = 75000
$ws.Range("M72").Value = -96236.00399999999
$ws.Range("N72").Value = -82488
$ws.Range("H94").Value = 2692.7144
$ws.Range("I94").Value = 3337.3333
$ws.Range("J94").Value = 2516.9092
$ws.Range("K94").Value = 3337.3333
$ws.Range("L94").Value = 2516.9092
$ws.Range("M94").Value = -2886.3333
$ws.Range("N94").Value = -3418.9092
$ws.Range("H107").Value = 371.58334
$ws.Range("I107").Value = 289.5
$ws.Range("K107").Value = 289.5
$ws.Range("M107").Value = 1630.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 186
$ws.Range("I8").Value = 186
$ws.Range("K8").Value = 558
$ws.Range("M8").Value = -419

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19999.834
$ws.Range("I70").Value = 19999.834
$ws.Range("K70").Value = 19999.834
$ws.Range("M70").Value = -19729.834
$ws.Range("H73").Value = 19999.834
$ws.Range("I73").Value = 19999.834
$ws.Range("K73").Value = 19999.834
$ws.Range("M73").Value = -19063.834
$ws.Range("H113").Value = 8476.799999999999
$ws.Range("I113").Value = 6095
$ws.Range("J113").Value = 10064.667
$ws.Range("K113").Value = 6095
$ws.Range("L113").Value = 10064.667
$ws.Range("M113").Value = -3925
$ws.Range("N113").Value = -14404.667
$ws.Range("H126").Value = 2860.625
$ws.Range("I126").Value = 2860.625
$ws.Range("K126").Value = 8581.875
$ws.Range("M126").Value = -6111.875
$ws.Range("H132").Value = 3258.8
$ws.Range("I132").Value = 2823.75
$ws.Range("K132").Value = 8471.25
$ws.Range("M132").Value = -5941.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7271.75
$ws.Range("I40").Value = 6452.2
$ws.Range("K40").Value = 6452.2
$ws.Range("M40").Value = -6316.2
$ws.Range("H55").Value = 1813.3158
$ws.Range("I55").Value = 182.1
$ws.Range("J55").Value = 3625.7778
$ws.Range("K55").Value = 182.1
$ws.Range("L55").Value = 3625.7778
$ws.Range("M55").Value = -9.099999999999994
$ws.Range("N55").Value = -3971.7778
$ws.Range("H61").Value = 3630.3928
$ws.Range("I61").Value = 2605.5625
$ws.Range("K61").Value = 2605.5625
$ws.Range("M61").Value = -2403.5625
$ws.Range("H93").Value = 3176.0557
$ws.Range("J93").Value = 4208
$ws.Range("L93").Value = 4208
$ws.Range("N93").Value = -6704
$ws.Range("H96").Value = 38833.168
$ws.Range("J96").Value = 38833.168
$ws.Range("L96").Value = 38833.168
$ws.Range("N96").Value = -44325.168
$ws.Range("H100").Value = 6716.5557
$ws.Range("I100").Value = 4767
$ws.Range("J100").Value = 8666.111000000001
$ws.Range("K100").Value = 4767
$ws.Range("L100").Value = 8666.111000000001
$ws.Range("M100").Value = -4226
$ws.Range("N100").Value = -9748.111000000001
$ws.Range("H113").Value = 3630.3928
$ws.Range("I113").Value = 2605.5625
$ws.Range("K113").Value = 2605.5625
$ws.Range("M113").Value = -435.5625
$ws.Range("H122").Value = 6151
$ws.Range("I122").Value = 5991.524
$ws.Range("K122").Value = 17974.572
$ws.Range("M122").Value = -15524.572

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14117.667
$ws.Range("J62").Value = 14238.229
$ws.Range("L62").Value = 14238.229
$ws.Range("N62").Value = -15486.229
$ws.Range("H65").Value = 14117.667
$ws.Range("J65").Value = 14238.229
$ws.Range("N65").Value = -77431.14499999999

Write-Output "Applied 221 cell updates across 8 sheets"